# Adds a "TOTAL" summary block (rows 27-29) to the salary statement sheet,
# mirroring the structure of the existing per-employee blocks (rows 15-17,
# 18-20, 21-23, 24-26) and summing across them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Merge the header-style cells on row 27 (same layout as rows 15/18/21/24)
#    Do this BEFORE copying formats over, so the merge operation doesn't
#    clobber the per-cell styles we are about to paste in.
# ---------------------------------------------------------------------
$ws.Range("B27:C27").Merge()
$ws.Range("G27:H27").Merge()

# ---------------------------------------------------------------------
# 2) Row 28 loses its B/C/D cells entirely in the target layout (unlike the
#    analogous row 25 which still has a B cell). Clear them out first so
#    PasteSpecial below won't leave stray cells behind.
# ---------------------------------------------------------------------
$ws.Range("B28:D28").Clear()

# ---------------------------------------------------------------------
# 3) Copy cell formatting (styles) from the last employee block (rows
#    24-26) onto the new TOTAL block (rows 27-29), fixing up the handful
#    of cells whose styles differ from a straight row copy.
# ---------------------------------------------------------------------

# Row 27 <- Row 24 (identical style layout)
$ws.Range("A24:U24").Copy()
$ws.Range("A27:U27").PasteSpecial(-4122)

# Row 28 <- Row 25 (identical except column E, which needs style 40 instead
# of "no style"; grab that from another style-40 cell such as L25)
$ws.Range("A25:A25").Copy()
$ws.Range("A28:A28").PasteSpecial(-4122)
$ws.Range("L25:L25").Copy()
$ws.Range("E28:E28").PasteSpecial(-4122)
$ws.Range("F25:U25").Copy()
$ws.Range("F28:U28").PasteSpecial(-4122)

# Row 29 <- Row 26 (identical except C29 needs style 44 instead of 45, and
# D29/E29 need style 46 instead of 44; source those from other cells with
# the right styles)
$ws.Range("A26:B26").Copy()
$ws.Range("A29:B29").PasteSpecial(-4122)
$ws.Range("B26:B26").Copy()
$ws.Range("C29:C29").PasteSpecial(-4122)
$ws.Range("M26:M26").Copy()
$ws.Range("D29:D29").PasteSpecial(-4122)
$ws.Range("M26:M26").Copy()
$ws.Range("E29:E29").PasteSpecial(-4122)
$ws.Range("F26:U26").Copy()
$ws.Range("F29:U29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Fill in the actual content: label, formulas and zero totals.
# ---------------------------------------------------------------------

# Row 27 - "TOTAL" header row
$ws.Cells.Item(27, 2).Value = "TOTAL"
$ws.Cells.Item(27, 5).Formula = "=E15+E18+E21+E24"
$ws.Cells.Item(27, 6).Value = 8
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Range("L27:Q27").Value = 0
$ws.Cells.Item(27, 18).Value = 0
$ws.Range("S27:T27").Value = 0

# Row 28
$ws.Cells.Item(28, 5).Formula = "=E16+E19+E22+E25"
$ws.Range("F28:K28").Value = 0
$ws.Range("L28:T28").Value = 0

# Row 29
$ws.Cells.Item(29, 4).Formula = "=D17+D20+D23+D26"
$ws.Cells.Item(29, 5).Formula = "=E17+E20+E23+E26"
$ws.Range("F29:I29").Value = 0
$ws.Range("K29:L29").Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Range("P29:S29").Value = 0
$ws.Cells.Item(29, 20).Value = 0

# ---------------------------------------------------------------------
# 5) Register the new shared string used for the label.
# ---------------------------------------------------------------------
# (Implicitly created by assigning the "TOTAL" text above.)
